# The sheet's data rows (2-17) are being reordered (a pure row permutation -
# every row's full content moves to a different row, nothing is added,
# removed or edited cell-by-cell). Build the new-row -> old-row mapping and
# apply it by reading every source row into memory first (so overwrites
# never clobber data we still need), then writing the rows back out in
# their new positions.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$firstRow = 2
$lastRow = 17
$firstColLetter = "A"
$lastColLetter = "AY"

# New row number -> old row number (where the data for the new row used to live).
$mapping = @{
    2  = 13
    3  = 2
    4  = 3
    5  = 4
    6  = 14
    7  = 5
    8  = 6
    9  = 7
    10 = 15
    11 = 8
    12 = 9
    13 = 10
    14 = 16
    15 = 17
    16 = 11
    17 = 12
}

# Most of the text columns hold plainly non-numeric text, which Excel's
# COM layer round-trips as text without any help. A few columns, though,
# contain text that LOOKS like a number or a date ("1" in I, "2012-09-26"
# in Y/AA) - left alone, assigning those through .Value2 would get
# reinterpreted as a real number/date. Force text formatting on just those
# columns before writing so the original string values survive untouched.
$textLikeNumberCols = @("I", "Y", "AA")
foreach ($col in $textLikeNumberCols) {
    $ws.Range("$col$firstRow`:$col$lastRow").NumberFormat = "@"
}

# Snapshot every source row's full data (A:AY) before any writes happen.
$rowsData = @{}
for ($r = $firstRow; $r -le $lastRow; $r++) {
    $rowsData[$r] = $ws.Range("$firstColLetter$r`:$lastColLetter$r").Value2
}

# Write each row back to its new position using the snapshot.
foreach ($newRow in $mapping.Keys) {
    $oldRow = $mapping[$newRow]
    $ws.Range("$firstColLetter$newRow`:$lastColLetter$newRow").Value2 = $rowsData[$oldRow]
}
